$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.938.19"
$ws.Range("E2").Value = "  +5.49%  "
$ws.Range("D3").Value = "2.418.87"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  +2.24%  "
$ws.Range("D9").Value = "2.449.96"
$ws.Range("E9").Value = "  +3.59%  "
$ws.Range("E10").Value = "  +6.17%  "
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("E12").Value = "  +3.25%  "
$ws.Range("E13").Value = "  +5.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000178"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.94%  "
$ws.Range("D16").Value = "2.971.81"
$ws.Range("E16").Value = "  +6.18%  "
$ws.Range("D17").Value = "62.676.81"
$ws.Range("E17").Value = "  +5.08%  "
$ws.Range("D18").Value = "2.458.91"
$ws.Range("E18").Value = "  +3.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("E22").Value = "  +2.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "621.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.13%  "
$ws.Range("D28").Value = "0.0₃0982"
$ws.Range("E28").Value = "  +7.48%  "
$ws.Range("D29").Value = "2.565.38"
$ws.Range("E29").Value = "  +3.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.24%  "
$ws.Range("E31").Value = "  +8.84%  "
$ws.Range("E32").Value = "  +4.63%  "
$ws.Range("E33").Value = "  +5.96%  "
$ws.Range("E34").Value = "  +5.20%  "
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.23%  "
$ws.Range("E37").Value = "  +2.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "152.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.13%  "
$ws.Range("E41").Value = "  +16.29%  "
$ws.Range("E42").Value = "  +7.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.45%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  -4.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "144.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.97%  "
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.601"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0517"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.61%  "
$ws.Range("E51").Value = "  +3.09%  "
